$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting the existing rows 34-78 down to 35-79
$ws.Rows(34).Insert()

# Populate the newly inserted row 34 with the new market record
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44579
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = 100112026
$ws.Range("G34").Value = "Haba"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 24000
$ws.Range("L34").Value = 24000
$ws.Range("M34").Value = 24000
$ws.Range("N34").Value = "`$/saco 25 kilos"
$ws.Range("O34").Value = "Región de La Araucanía"
$ws.Range("P34").Value = 960
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
